# Add stats bars for enemies and connect them to components.
# - mark "add enemy hp bar and ? Mp" (row 16) as done, with a done-at date
# - mark "player attacks reduce enemy hp (hp shown)" (row 18) as done
# - extend the matching dev-log entry on the Logs sheet with the new work
# - move the active selection/tab to the Logs sheet

$wb = $excel.ActiveWorkbook

# --- Sheet "TODO Before 0.0.1" -------------------------------------------
$ws1 = $wb.Worksheets.Item("TODO Before 0.0.1")

# Row 16: "add enemy hp bar and ? Mp" -> done, done-at 6/30/2024 (45473)
$ws1.Range("C16").Value = "done"
$d16 = $ws1.Range("D16")
$d16.ClearFormats()
$d16.Value = 45473
$d16.NumberFormat = "mm-dd-yy"

# Row 18: "player attacks reduce enemy hp (hp shown)" -> done (date already set)
$ws1.Range("C18").Value = "done"

# Move the sheet's saved selection off the old spot (no longer the active tab)
$ws1.Range("C14").Select()

# --- Sheet "Logs" ----------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Logs")

# Append the new work to the existing dev-log line about enemy refactor
$ws3.Range("B41").Value = "finish refactoring enemy, fix bug when player shoot underground while camera close to wall, add stat bars for enemies"

# "Logs" becomes the active sheet/tab, with a new selection
$ws3.Activate()
$ws3.Range("B42").Select()
